$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with plain (lower-case, no-prefix) labels that
# correspond to the "iaest-measure:*" row that used to be at row 2.
$ws.Range("A2").Value = "municipio"
$ws.Range("B2").Value = "entidad-singular"
$ws.Range("C2").Value = "personas"
$ws.Range("D2").Value = "nucleo"
$ws.Range("E2").Value = "comarca"
$ws.Range("F2").Value = "provincia"
$ws.Range("G2").Value = "nucleodiseminado"

# Match style of the rest of the sheet (style index 1 / font "Arial 10").
$ws.Range("A1:G1").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
